$d = $word.ActiveDocument

# Replace the text of the first paragraph (the "<header>" line) with the new <img> tag line.
$d.Paragraphs.Item(1).Range.Text = '<img src="profile.jpg" alt="Profile photo of John Doe">'

# Remove all the remaining paragraphs (everything that followed the first paragraph),
# i.e. delete from the end of the (now updated) first paragraph through the end of the document.
$start = $d.Paragraphs.Item(1).Range.End
$end = $d.Content.End
$r = $d.Range($start, $end)
$r.Delete()
